$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.655.81"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "3.524.73"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'613.00"
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("D6").Value = "'173.97"
$ws.Range("E6").Value = "  +1.05%  "
$ws.Range("D7").Value = "'0.610"
$ws.Range("E7").Value = "  -1.40%  "
$ws.Range("D8").Value = "3.517.29"
$ws.Range("E8").Value = "  -0.96%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "'0.197"
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("D11").Value = "'7.17"
$ws.Range("E11").Value = "  +3.82%  "
$ws.Range("D12").Value = "'0.591"
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("D13").Value = "'46.62"
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("D14").Value = "'0.0000277"
$ws.Range("E14").Value = "  -0.40%  "
$ws.Range("D15").Value = "4.092.31"
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("D16").Value = "'8.47"
$ws.Range("E16").Value = "  +1.02%  "
$ws.Range("D17").Value = "'615.87"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "70.618.95"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.511.82"
$ws.Range("E19").Value = "  -1.11%  "
$ws.Range("E20").Value = "  +1.98%  "
$ws.Range("D21").Value = "'17.77"
$ws.Range("E21").Value = "  +2.20%  "
$ws.Range("D22").Value = "'0.886"
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("D23").Value = "'9.02"
$ws.Range("E23").Value = "  -4.23%  "
$ws.Range("D24").Value = "'98.84"
$ws.Range("E24").Value = "  +2.08%  "
$ws.Range("D25").Value = "'15.73"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").Value = "'3.78"
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").Value = "'2.60"
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("D29").Value = "'33.96"
$ws.Range("E29").Value = "  +1.04%  "
$ws.Range("D30").Value = "'9.23"
$ws.Range("E30").Value = "  +1.95%  "
$ws.Range("D31").Value = "'8.17"
$ws.Range("E31").Value = "  -3.63%  "
$ws.Range("D32").Value = "'3.02"
$ws.Range("E32").Value = "  -1.78%  "
$ws.Range("E33").Value = "  -0.73%  "
$ws.Range("D34").Value = "'6.88"
$ws.Range("E34").Value = "  -1.47%  "
$ws.Range("D35").Value = "'617.93"
$ws.Range("E35").Value = "  +7.52%  "
$ws.Range("D36").Value = "'0.101"
$ws.Range("E36").Value = "  -0.48%  "
$ws.Range("D37").Value = "'0.0494"
$ws.Range("E37").Value = "  +5.38%  "
$ws.Range("D38").Value = "'10.85"
$ws.Range("E38").Value = "  +0.40%  "
$ws.Range("D39").Value = "'3.52"
$ws.Range("E39").Value = "  -3.09%  "
$ws.Range("D40").Value = "'57.00"
$ws.Range("E40").Value = "  -0.98%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").Value = "'0.145"
$ws.Range("E42").Value = "  +1.62%  "
$ws.Range("B43").Value = "PEPE"
$ws.Range("C43").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D43").Value = "0.0₃0744"
$ws.Range("E43").Value = "  +5.55%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "3.378.39"
$ws.Range("E44").Value = "  +0.82%  "
$ws.Range("D45").Value = "'0.312"
$ws.Range("E45").Value = "  -2.78%  "
$ws.Range("D46").Value = "'2.94"
$ws.Range("E46").Value = "  -2.26%  "
$ws.Range("D47").Value = "'32.33"
$ws.Range("E47").Value = "  -2.03%  "
$ws.Range("D48").Value = "'2.58"
$ws.Range("E48").Value = "  -2.11%  "
$ws.Range("D49").Value = "'0.131"
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("D50").Value = "'133.74"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("E51").Value = "  -0.01%  "
